# Weekly price update: insert a new Primera/Segunda pair of rows for
# "Betarraga" (Vega Monumental Concepción) dated 2022-07-12 (serial 44754)
# at the top of the data block (row 270), pushing the existing history
# down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 270-271; everything from old row 270 onward
# shifts down to 272 onward (same as Excel's row-insert behaviour).
$ws.Range("A270:A271").EntireRow.Insert()

# New row 270 - "Primera" quality entry for the week of 2022-07-12.
$ws.Range("A270").Value2 = 11
$ws.Range("B270").Value2 = "Vega Monumental Concepción"
$ws.Range("C270").Value2 = "Bíobío"
$ws.Range("D270").Value2 = 44754
$ws.Range("E270").Value2 = 8
$ws.Range("F270").Value2 = 100114014
$ws.Range("G270").Value2 = "Betarraga"
$ws.Range("H270").Value2 = "Sin especificar"
$ws.Range("I270").Value2 = "Primera"
$ws.Range("J270").Value2 = 800
$ws.Range("K270").Value2 = 600
$ws.Range("L270").Value2 = 700
$ws.Range("M270").Value2 = 650
$ws.Range("N270").Value2 = "`$/paquete 5 unidades"
$ws.Range("O270").Value2 = "Región Metropolitana"
$ws.Range("P270").Value2 = 130
$ws.Range("Q270").Value2 = 5
$ws.Range("R270").Value2 = "Hortaliza"

# New row 271 - "Segunda" quality entry for the same week.
$ws.Range("A271").Value2 = 11
$ws.Range("B271").Value2 = "Vega Monumental Concepción"
$ws.Range("C271").Value2 = "Bíobío"
$ws.Range("D271").Value2 = 44754
$ws.Range("E271").Value2 = 8
$ws.Range("F271").Value2 = 100114014
$ws.Range("G271").Value2 = "Betarraga"
$ws.Range("H271").Value2 = "Sin especificar"
$ws.Range("I271").Value2 = "Segunda"
$ws.Range("J271").Value2 = 400
$ws.Range("K271").Value2 = 500
$ws.Range("L271").Value2 = 500
$ws.Range("M271").Value2 = 500
$ws.Range("N271").Value2 = "`$/paquete 5 unidades"
$ws.Range("O271").Value2 = "Región Metropolitana"
$ws.Range("P271").Value2 = 100
$ws.Range("Q271").Value2 = 5
$ws.Range("R271").Value2 = "Hortaliza"

# Make sure the date cells keep the same number format as the rest of
# column D (the insert already carries style down, but set explicitly
# to be safe).
$ws.Range("D270:D271").NumberFormat = $ws.Range("D272").NumberFormat
